$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.828.82"
$ws.Range("E2").Value = "  -1.04%  "

$ws.Range("D3").Value = "2.341.93"
$ws.Range("E3").Value = "  +1.17%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.18"
$ws.Range("E5").Value = "  -1.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.46"
$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.511"
$ws.Range("E7").Value = "  -4.60%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  -2.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.13"
$ws.Range("E10").Value = "  -2.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.46"
$ws.Range("E11").Value = "  +0.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0797"
$ws.Range("E12").Value = "  -2.38%  "

$ws.Range("E13").Value = "  -0.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.84"
$ws.Range("E14").Value = "  -3.51%  "

$ws.Range("D15").Value = "2.711.96"
$ws.Range("E15").Value = "  +1.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.44"
$ws.Range("E16").Value = "  +2.94%  "

$ws.Range("D17").Value = "2.346.70"
$ws.Range("E17").Value = "  +1.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.799"
$ws.Range("E18").Value = "  -1.71%  "

$ws.Range("D19").Value = "42.785.31"
$ws.Range("E19").Value = "  -0.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.25"
$ws.Range("E20").Value = "  +1.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.71"
$ws.Range("E21").Value = "  -6.86%  "

$ws.Range("D22").Value = "0.0₃0905"
$ws.Range("E22").Value = "  -1.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.40"
$ws.Range("E23").Value = "  -1.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "237.02"
$ws.Range("E24").Value = "  -1.88%  "

$ws.Range("E25").Value = "  -1.44%  "

$ws.Range("E26").Value = "  -2.59%  "

$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.32"
$ws.Range("E28").Value = "  +2.40%  "

$ws.Range("E29").Value = "  -3.86%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.24"
$ws.Range("E30").Value = "  +5.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.43"
$ws.Range("E31").Value = "  -5.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.40"
$ws.Range("E32").Value = "  -2.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "159.97"
$ws.Range("E33").Value = "  -4.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("E35").Value = "  -3.55%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.84"
$ws.Range("E36").Value = "  -0.79%  "

$ws.Range("E37").Value = "  +3.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0728"
$ws.Range("E38").Value = "  -2.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.59"
$ws.Range("E39").Value = "  +6.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.02"
$ws.Range("E40").Value = "  -4.22%  "

$ws.Range("E41").Value = "  +1.60%  "

$ws.Range("E42").Value = "  -3.62%  "

$ws.Range("E43").Value = "  -3.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.51"
$ws.Range("E44").Value = "  +8.59%  "

$ws.Range("D45").Value = "2.022.57"
$ws.Range("E45").Value = "  +2.42%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0286"
$ws.Range("E46").Value = "  -1.69%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.04"
$ws.Range("E47").Value = "  -4.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.47"
$ws.Range("E48").Value = "  +6.66%  "

$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "56.94"
$ws.Range("E50").Value = "  +2.01%  "

$ws.Range("E51").Value = "  -2.26%  "
